$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Add new worksheet "PreguntasFrecuentes" right after "Empleados"
# ------------------------------------------------------------------
$empleados = $wb.Worksheets.Item("Empleados")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $empleados)
$ws.Name = "PreguntasFrecuentes"

# ------------------------------------------------------------------
# 2) Column widths (closest achievable match to the bestFit widths of
#    the target file, compensating for this engine's pixel rounding)
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 47.333333333333336
$ws.Columns.Item(2).ColumnWidth = 15.333333333333334
$ws.Columns.Item(3).ColumnWidth = 8.666666666666666
$ws.Columns.Item(4).ColumnWidth = 47.166666666666664
$ws.Columns.Item(5).ColumnWidth = 17.333333333333332
$ws.Columns.Item(6).ColumnWidth = 7.5
$ws.Columns.Item(7).ColumnWidth = 18.5
$ws.Columns.Item(8).ColumnWidth = 12.666666666666666
$ws.Columns.Item(9).ColumnWidth = 30.833333333333332
$ws.Columns.Item(10).ColumnWidth = 12.666666666666666

# ------------------------------------------------------------------
# 3) Header row
# ------------------------------------------------------------------
$ws.Range("A1").Value = "CasoPrueba"
$ws.Range("B1").Value = "Usuario"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Url"
$ws.Range("E1").Value = "ResultadoEsperado"
$ws.Range("F1").Value = "Browser"
$ws.Range("G1").Value = "Categoria"
$ws.Range("H1").Value = "Activar"
$ws.Range("I1").Value = "Pregunta"
$ws.Range("J1").Value = "Respuesta"
$ws.Range("K1").Value = "TipoPrueba"

$url = "https://inter-wepid-dev.azurewebsites.net/Admin/"

# ------------------------------------------------------------------
# 4) Data rows 2-9
#    NOTE: "true"/"false" are entered as plain text (leading apostrophe)
#    so that they become shared-string text cells, matching the source
#    sheets, instead of native Excel booleans.
# ------------------------------------------------------------------
$data = @(
    @("Alta_Exitosa_PreguntaCamposObligatorios", "admin@inter.mx", 12345678, $url, "'true", "Chrome", "Servicio / Producto", "'true",  "Tiene dudas sobre el producto", "xxxx",     "Crea"),
    @("Alta_Exitosa_PreguntaCamposObligatorios", "admin@inter.mx", 12345678, $url, "'true", "Chrome", "Robo de identidad",    "'true",  "Como proteger su identidad",    "xxxx",     "Crea"),
    @("Alta_Exitosa_PreguntaCamposObligatorios", "admin@inter.mx", 12345678, $url, "'true", "Chrome", "Contacto",             "'true",  "Desea contactarnos",            "xxxx",     "Crea"),
    @("Alta_Exitosa_PreguntaCamposObligatorios", "admin@inter.mx", 12345678, $url, "'true", "Chrome", "Quienes somos",        "'false", "Informacion de WEPID",          "xxxx",     "Crea"),
    @("Editar_Pregunta",                         "admin@inter.mx", 12345678, $url, "'true", "Chrome", "Servicio / Producto", "'true",  "La pregunta sera dada de baja", "xxxxxxxx", "editar"),
    @("Eliminar_Pregunta",                       "admin@inter.mx", 12345678, $url, "'true", "Chrome", "Robo de identidad",    "'true",  "N/A",                           $null,      "Eliminar"),
    @("Eliminar_Pregunta",                       "admin@inter.mx", 12345678, $url, "'true", "Chrome", "Contacto",             "'true",  "N/A",                           $null,      "Eliminar"),
    @("Eliminar_Pregunta",                       "admin@inter.mx", 12345678, $url, "'true", "Chrome", "Quienes somos",        "'true",  "N/A",                           $null,      "Eliminar")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $ws.Cells.Item($rowIndex, 9).Value = $row[8]
    if ($null -ne $row[9]) {
        $ws.Cells.Item($rowIndex, 10).Value = $row[9]
    }
    $ws.Cells.Item($rowIndex, 11).Value = $row[10]
    $ws.Hyperlinks.Add($ws.Cells.Item($rowIndex, 4), $url) | Out-Null
    $rowIndex++
}

# ------------------------------------------------------------------
# 5) Styles - copy direct formatting from equivalent cells on the
#    "Empleados" sheet, which already contains each needed cell style.
# ------------------------------------------------------------------
$empleados.Range("A1").Copy() | Out-Null
$ws.Range("A1:K1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats : header style (s=1)

$empleados.Range("A2").Copy() | Out-Null
$ws.Range("A2:F9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats : plain bordered style (s=2)
$ws.Range("I2:K9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats : plain bordered style (s=2)

$empleados.Range("B2").Copy() | Out-Null
$ws.Range("B2:B9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats : hyperlink-font bordered style (s=3)
$ws.Range("D2:D9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats : hyperlink-font bordered style (s=3)

$empleados.Range("G2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null
$ws.Range("G9").PasteSpecial(-4122) | Out-Null

$empleados.Range("G4").Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4122) | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null

$empleados.Range("H2").Copy() | Out-Null
$ws.Range("H2:H9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats : currency-hyperlink bordered style (s=6)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 6) Selection and active sheet
#    Select A2 on "Empleados" first (it is no longer the active tab),
#    then select A9 on the new sheet and make it the active/visible tab.
# ------------------------------------------------------------------
$empleados.Range("A2").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("A9").Select() | Out-Null
